# Ravirala_LabExam03Grading.xlsx - grading update for question 12
# (CustomerMapping Class section): the grader revised the deduction
# comment for row 20 (S.No 12) from a generic "(-5) For incorrect logic."
# note to two itemised (-1) deductions, bumped the awarded points for
# that question from 5 to 8, widened column F and wrapped/grew the row
# so the longer, two-line comment is fully visible.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Points awarded for question 12 (S. No column) changed 5 -> 8.
$ws.Range("E20").Value = 8

# Replace the grading comment in F20 with the new, two-line deduction text.
$ws.Range("F20").Value = "(-1) For incorrect creation of linked list." + [char]10 + "(-1) For not adding new customer into hash map."

# The new comment needs to wrap within the cell instead of overflowing.
$ws.Range("F20").WrapText = $true

# Widen column F so the long comment text fits, and grow row 20 to show
# both wrapped lines.
$ws.Columns.Item(6).ColumnWidth = 129.42578125
$ws.Rows.Item(20).RowHeight = 30

# The totals below (Total row E26 and Grand Total E38) are driven by
# SUM() formulas over this range, so they recalculate automatically.
